# Applies the "updated first two scenarios" change:
#  - Create sheet: add STATUS_CODE / createdAt columns (D, E) with values
#    201 / 2024 (as text) for every data row, sets text number format on
#    those columns, and updates the active selection.
#  - Update-PUT sheet: updates the active selection only.

$wb = $excel.ActiveWorkbook
$wsCreate = $wb.Worksheets.Item("Create")
$wsUpdatePut = $wb.Worksheets.Item("Update-PUT")

# --- Create sheet: new STATUS_CODE / createdAt columns -----------------
# Apply the text number format *before* assigning values so "201"/"2024"
# are stored as text (shared strings) rather than being coerced to numbers.
$wsCreate.Range("D1:E5").NumberFormat = "@"

$wsCreate.Range("D1").Value = "STATUS_CODE"
$wsCreate.Range("E1").Value = "createdAt"

$wsCreate.Range("D2").Value = "201"
$wsCreate.Range("D3").Value = "201"
$wsCreate.Range("D4").Value = "201"
$wsCreate.Range("D5").Value = "201"

$wsCreate.Range("E2").Value = "2024"
$wsCreate.Range("E3").Value = "2024"
$wsCreate.Range("E4").Value = "2024"
$wsCreate.Range("E5").Value = "2024"

$wsCreate.Columns.Item(4).ColumnWidth = 12
$wsCreate.Columns.Item(5).ColumnWidth = 8

# --- Active selections ---------------------------------------------------
# Update-PUT becomes inactive, its saved selection moves to A17.
$wsUpdatePut.Range("A17").Select() | Out-Null

# Create becomes the tab-selected / active sheet, selection moves to D8.
$wsCreate.Select()
$wsCreate.Range("D8").Select() | Out-Null
